# Updated symbol list on Fri Dec 23 11:35:25 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates. These values must remain plain text (as they were
# stored as inline/shared strings, not numbers), so we temporarily force a
# text number-format before assigning, then clear the format again so no
# stray style id is left behind on the cell.
$priceUpdates = @{
    "D2"  = "246.11"
    "D3"  = "22.01"
    "D4"  = "5.428"
    "D5"  = "0.05836"
    "D7"  = "6.348"
    "D8"  = "0.8072"
    "D9"  = "0.9665"
    "D11" = "0.07431"
    "D12" = "0.03326"
    "D13" = "0.03035"
    "D14" = "4.175"
    "D15" = "0.09404"
    "D16" = "0.001597"
    "D17" = "0.04814"
    "D19" = "0.006131"
    "D20" = "0.004106"
    "D21" = "0.0009975"
    "D23" = "3.702"
    "D24" = "2.214"
    "D25" = "0.3209"
    "D26" = "0.1296"
    "D40" = "0.03874"
    "D41" = "0.006662"
    "D42" = "0.1075"
    "D43" = "0.003001"
    "D44" = "0.006687"
    "D45" = "0.00005612"
    "D47" = "0.4201"
    "D48" = "0.1457"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.ClearFormats()
}

# Volume(1h) label column (E) updates - these are plain text already, so a
# direct value assignment keeps them as text.
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E27").Value = "26UpBotsUBXTWorstin24h"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"
